# Update "想去人数" (interest count) values in column F on the
# "展览" and "全部类型" worksheets, per the latest data refresh.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value
$exhibitionUpdates = @{
    "F3"  = 1112
    "F4"  = 1920
    "F6"  = 1252
    "F8"  = 38
    "F13" = 815
    "F14" = 242
    "F21" = 74
    "F23" = 188
    "F32" = 423
}

foreach ($cellRef in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range($cellRef).Value = $exhibitionUpdates[$cellRef]
}

# Sheet "全部类型": row -> new value
$allTypesUpdates = @{
    "F4"  = 1112
    "F5"  = 1920
    "F7"  = 1252
    "F10" = 38
    "F15" = 815
    "F16" = 242
    "F29" = 74
    "F31" = 188
    "F46" = 423
}

foreach ($cellRef in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range($cellRef).Value = $allTypesUpdates[$cellRef]
}
